$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cell from "Responses" to "Response"
$ws.Range("A1").Value = "Response"

# The diff drops the stale view state (topLeftCell="A11" and the
# G23 selection), i.e. the saved view goes back to showing/selecting A1.
$ws.Range("A1").Select()
